# Update the "取得日時" (retrieved timestamp) column for all data rows
# on the active sheet ("ランサーズ") from 2025-11-20 01:18:39 to
# 2025-11-20 01:47:58, reflecting a new scrape/append run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-11-20 01:47:58"

$ws.Range("A2:A15").Value = $newTimestamp
